$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.350.41"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.004.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +7.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9969"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8028"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +69.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "259.84"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9968"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3623"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +26.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.68"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +32.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07070"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8449"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +17.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08090"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.000.94"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "101.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.640"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.320.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +13.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.925"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +12.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007941"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.263.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9955"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9969"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.200"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +15.11%  "

$ws.Range("E25").Value = "  +13.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1540"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +59.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("E28").Value = "  +7.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.418"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +28.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.635"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.646"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +10.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.355"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.408"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05224"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +9.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.226"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7676"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.804"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.48%  "

$ws.Range("E38").Value = "  +6.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.954"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.04"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.25%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.684"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.37%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.206"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +14.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4746"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8603"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9975"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.672"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.95%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.949"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4380"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +12.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1193"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +15.33%  "
